$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = "31/12/2006"; C = 73.4913550552766 },
    @{ Row = 3;  B = "31/12/2010"; C = 74.3604259197093 },
    @{ Row = 4;  B = "31/12/2014"; C = 75.1981752438586 },
    @{ Row = 5;  B = "31/12/2018"; C = 75.9902268405237 },
    @{ Row = 6;  B = "31/12/2022"; C = 75.44804410408101 },
    @{ Row = 7;  B = "31/12/2024"; C = 76.60571098458711 },
    @{ Row = 8;  B = "31/12/2006"; C = 73.1680763106949 },
    @{ Row = 9;  B = "31/12/2010"; C = 74.17092673285011 },
    @{ Row = 10; B = "31/12/2014"; C = 74.7570772999095 },
    @{ Row = 11; B = "31/12/2018"; C = 75.4522008450428 },
    @{ Row = 12; B = "31/12/2022"; C = 74.9193910754588 },
    @{ Row = 13; B = "31/12/2024"; C = 76.1510992616663 },
    @{ Row = 14; B = "31/12/2006"; C = 73.2226792128026 },
    @{ Row = 15; B = "31/12/2010"; C = 74.09392402191629 },
    @{ Row = 16; B = "31/12/2014"; C = 74.1091720876165 },
    @{ Row = 17; B = "31/12/2018"; C = 75.040398524986 },
    @{ Row = 18; B = "31/12/2022"; C = 75.26526792012019 },
    @{ Row = 19; B = "31/12/2024"; C = 76.3617434255048 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
